# Updated CHE_grids model - 2025-08-21 14:14
#
# 1) Extend the shared "value" list in column T (AVA sheet) from
#    "e_demand,ev_battery" to
#    "e_demand,ev_battery,H2prd_Elc_PEM,H2prd_Elc_ALK"
#    for every row that currently holds that text.
# 2) Widen column T to match the other "value" columns (E, J) and
#    drop the bestFit auto-size flag now that it holds a custom width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVA")

$oldText = "e_demand,ev_battery"
$newText = "e_demand,ev_battery,H2prd_Elc_PEM,H2prd_Elc_ALK"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 20)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}

# Match the width used by the analogous columns E and J, and clear bestFit.
$ws.Columns.Item(20).ColumnWidth = 29.8
